# Appends 9 new run-summary rows (212-220) to the "Run Info" sheet, matching
# the pattern of the existing rows (same columns/format), and updates the
# sheet's selection to the new last row, as produced by entering the data
# directly at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Run Info")
$ws.Activate()

# Each inner array is one new row: A..P (Date, RCP, E, everyx, workers,
# run time, 85-2010 bleaching, sBleach 1, cBleach 1, sRecSeedMult 1,
# cRecSeedMult 1, cSeedThreshMult 1, pMin, pMax, exponent, divisor).
$newRows = @(
    @(42963.427152777775, "rcp45", 1, 1, 6, 43.311756045473004, 5.0069930069930075, 0.3, 0.1, 4, 4, 2, 0.36, 1.5, 0.46, 4.8444000000000003),
    @(42963.472777777781, "rcp45", 1, 1, 6, 42.049631260625638, 5.0069930069930075, 0.3, 0.1, 4, 4, 2, 0.36, 1.5, 0.46, 4.8444000000000003),
    @(42963.473807870374, "rcp45", 1, 1, 11, 39.755159423341347, 5.0069930069930075, 0.3, 0.1, 4, 4, 2, 0.36, 1.5, 0.46, 4.8444000000000003),
    @(42963.474479166667, "rcp45", 1, 1, 6, 35.52102600006976, 5.0069930069930075, 0.3, 0.1, 4, 4, 2, 0.36, 1.5, 0.46, 4.8444000000000003),
    @(42963.692175925928, "rcp45", 1, 10000, 1, 14.547541572176568, 0, 0.3, 0.1, 4, 4, 2, 0.36, 1.5, 0.46, 4.8444000000000003),
    @(42963.696458333332, "rcp45", 1, 10000, 1, 13.107532958071516, 0, 0.3, 0.1, 4, 4, 2, 0.36, 1.5, 0.46, 4.8444000000000003),
    @(42963.711736111109, "rcp45", 1, 10000, 1, 14.766424637620409, 0, 0.3, 0.1, 4, 4, 2, 0.36, 1.5, 0.46, 4.8444000000000003),
    @(42963.712789351855, "rcp45", 1, 10000, 1, 13.574009145926844, 0, 0.3, 0.1, 4, 4, 2, 0.36, 1.5, 0.46, 4.8444000000000003),
    @(42963.743738425925, "rcp45", 1, 10000, 1, 14.438064040687246, 0, 0.3, 0.1, 4, 4, 2, 0.36, 1.5, 0.46, 4.8444000000000003)
)

$lastExistingRow = 211
$startRow = $lastExistingRow + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]

    # Carry the date/time number format from the row above into column A
    # so the new date cells keep the same style as the rest of the table.
    $ws.Cells.Item($r - 1, 1).Copy($ws.Cells.Item($r, 1))

    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}

$lastRow = $startRow + $newRows.Count - 1
$lastCol = 16

$ws.Range("A170").Select()
$ws.Application.ActiveWindow.ScrollRow = 170
$ws.Cells.Item($lastRow, 1).Activate()
$ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, $lastCol)).Select()
